# Scheduled data refresh: recompute market-price-derived columns (H:N)
# currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# for a batch of Leve rows across several Sheets (per the latest market snapshot).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row19 - Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 848
$ws.Range("I19").Value = 600
$ws.Range("J19").Value = 930.6667
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 930.6667
$ws.Range("M19").Value = -425
$ws.Range("N19").Value = -1280.6667

# ALC!row98 - The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 125000740
$ws.Range("I98").Value = 178571920
$ws.Range("J98").Value = 1333.3334
$ws.Range("K98").Value = 178571920
$ws.Range("L98").Value = 1333.3334
$ws.Range("M98").Value = -178570422
$ws.Range("N98").Value = -4329.3334

# ALC!row122 - Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 125000740
$ws.Range("I122").Value = 178571920
$ws.Range("J122").Value = 1333.3334
$ws.Range("K122").Value = 535715760
$ws.Range("L122").Value = 4000.0002
$ws.Range("M122").Value = -535713310
$ws.Range("N122").Value = -8900.0002

# ALC!row138 - All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1528.2354
$ws.Range("I138").Value = 876.89795
$ws.Range("J138").Value = 3208
$ws.Range("K138").Value = 2630.69385
$ws.Range("L138").Value = 9624
$ws.Range("M138").Value = 2509.30615
$ws.Range("N138").Value = -19904

$ws = $wb.Worksheets.Item("ARM")
# ARM!row45 - Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 15198011
$ws.Range("I45").Value = 20896784
$ws.Range("J45").Value = 1283.1666
$ws.Range("K45").Value = 20896784
$ws.Range("L45").Value = 1283.1666
$ws.Range("M45").Value = -20896407
$ws.Range("N45").Value = -2037.1666

$ws = $wb.Worksheets.Item("BSM")
# BSM!row86 - Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 970889
$ws.Range("I86").Value = 2160.3333
$ws.Range("J86").Value = 1939617.6
$ws.Range("K86").Value = 2160.3333
$ws.Range("L86").Value = 1939617.6
$ws.Range("M86").Value = -1037.3333
$ws.Range("N86").Value = -1941863.6

# BSM!row89 - Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 970889
$ws.Range("I89").Value = 2160.3333
$ws.Range("J89").Value = 1939617.6
$ws.Range("K89").Value = 10801.6665
$ws.Range("L89").Value = 9698088
$ws.Range("M89").Value = -5185.666499999999
$ws.Range("N89").Value = -9709320

$ws = $wb.Worksheets.Item("CRP")
# CRP!row31 - Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1351.2778
$ws.Range("I31").Value = 992.7045000000001
$ws.Range("J31").Value = 1914.75
$ws.Range("K31").Value = 992.7045000000001
$ws.Range("L31").Value = 1914.75
$ws.Range("M31").Value = -697.7045000000001
$ws.Range("N31").Value = -2504.75

# CRP!row34 - Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1351.2778
$ws.Range("I34").Value = 992.7045000000001
$ws.Range("J34").Value = 1914.75
$ws.Range("K34").Value = 992.7045000000001
$ws.Range("L34").Value = 1914.75
$ws.Range("M34").Value = -790.7045000000001
$ws.Range("N34").Value = -2318.75

# CRP!row58 - You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 11628990
$ws.Range("I58").Value = 14286729
$ws.Range("J58").Value = 1377.6875
$ws.Range("K58").Value = 14286729
$ws.Range("L58").Value = 1377.6875
$ws.Range("M58").Value = -14286526
$ws.Range("N58").Value = -1783.6875

# CRP!row99 - O Pine / Pine Lumber
$ws.Range("H99").Value = 142859140
$ws.Range("I99").Value = 250001340
$ws.Range("J99").Value = 2833.3333
$ws.Range("K99").Value = 250001340
$ws.Range("L99").Value = 2833.3333
$ws.Range("M99").Value = -249999842
$ws.Range("N99").Value = -5829.3333

# CRP!row122 - Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 22727982
$ws.Range("I122").Value = 35714910
$ws.Range("J122").Value = 849.75
$ws.Range("K122").Value = 107144730
$ws.Range("L122").Value = 2549.25
$ws.Range("M122").Value = -107142280
$ws.Range("N122").Value = -7449.25

# CRP!row126 - A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 142859140
$ws.Range("I126").Value = 250001340
$ws.Range("J126").Value = 2833.3333
$ws.Range("K126").Value = 750004020
$ws.Range("L126").Value = 8499.999899999999
$ws.Range("M126").Value = -750001550
$ws.Range("N126").Value = -13439.9999

# CRP!row132 - Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 11906421
$ws.Range("I132").Value = 1467.4286
$ws.Range("J132").Value = 23811374
$ws.Range("K132").Value = 4402.2858
$ws.Range("L132").Value = 71434122
$ws.Range("M132").Value = -1872.2858
$ws.Range("N132").Value = -71439182

# CRP!row134 - Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1815.36
$ws.Range("I134").Value = 1645.4546
$ws.Range("J134").Value = 1948.8572
$ws.Range("K134").Value = 4936.3638
$ws.Range("L134").Value = 5846.571599999999
$ws.Range("M134").Value = -2401.3638
$ws.Range("N134").Value = -10916.5716

# CRP!row136 - Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 11628990
$ws.Range("I136").Value = 14286729
$ws.Range("J136").Value = 1377.6875
$ws.Range("K136").Value = 42860187
$ws.Range("L136").Value = 4133.0625
$ws.Range("M136").Value = -42857637
$ws.Range("N136").Value = -9233.0625

$ws = $wb.Worksheets.Item("CUL")
# CUL!row92 - Oh No Udon / Gyr Abanian Flour
$ws.Range("H92").Value = 5767.857
$ws.Range("I92").Value = 230.5
$ws.Range("J92").Value = 7070.7646
$ws.Range("K92").Value = 691.5
$ws.Range("L92").Value = 21212.2938
$ws.Range("M92").Value = 556.5
$ws.Range("N92").Value = -23708.2938

$ws = $wb.Worksheets.Item("GSM")
# GSM!row70 - Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 4388.914
$ws.Range("I70").Value = 3794.6191
$ws.Range("J70").Value = 5280.357
$ws.Range("K70").Value = 3794.6191
$ws.Range("L70").Value = 5280.357
$ws.Range("M70").Value = -3524.6191
$ws.Range("N70").Value = -5820.357

# GSM!row73 - Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 4388.914
$ws.Range("I73").Value = 3794.6191
$ws.Range("J73").Value = 5280.357
$ws.Range("K73").Value = 3794.6191
$ws.Range("L73").Value = 5280.357
$ws.Range("M73").Value = -2858.6191
$ws.Range("N73").Value = -7152.357

# GSM!row113 - Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 1037
$ws.Range("I113").Value = 1037
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1037
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1133
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# LTW!row7 - Tan Before the Ban / Leather
$ws.Range("H7").Value = 1664.5834
$ws.Range("I7").Value = 1634.091
$ws.Range("K7").Value = 1634.091
$ws.Range("M7").Value = -1522.091

# LTW!row22 - Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 1802.1305
$ws.Range("I22").Value = 699.8333
$ws.Range("J22").Value = 2191.1765
$ws.Range("K22").Value = 699.8333
$ws.Range("L22").Value = 2191.1765
$ws.Range("M22").Value = -404.8333
$ws.Range("N22").Value = -2781.1765

# LTW!row27 - Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 1802.1305
$ws.Range("I27").Value = 699.8333
$ws.Range("J27").Value = 2191.1765
$ws.Range("K27").Value = 699.8333
$ws.Range("L27").Value = 2191.1765
$ws.Range("M27").Value = -592.8333
$ws.Range("N27").Value = -2405.1765

# LTW!row61 - Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1666.6666
$ws.Range("I61").Value = 1647.3529
$ws.Range("J61").Value = 1995
$ws.Range("K61").Value = 1647.3529
$ws.Range("L61").Value = 1995
$ws.Range("M61").Value = -1445.3529
$ws.Range("N61").Value = -2399

# LTW!row113 - Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1666.6666
$ws.Range("I113").Value = 1647.3529
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 1647.3529
$ws.Range("L113").Value = 1995
$ws.Range("M113").Value = 522.6470999999999
$ws.Range("N113").Value = -6335

# LTW!row126 - Battered Books / Saiga Leather
$ws.Range("H126").Value = 1664.5834
$ws.Range("I126").Value = 1634.091
$ws.Range("K126").Value = 4902.272999999999
$ws.Range("M126").Value = -2432.272999999999

# LTW!row132 - Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 15877988
$ws.Range("I132").Value = 23811078
$ws.Range("K132").Value = 71433234
$ws.Range("M132").Value = -71430704

# LTW!row136 - Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 20051710
$ws.Range("I136").Value = 3324080.2
$ws.Range("J136").Value = 71429430
$ws.Range("K136").Value = 9972240.600000001
$ws.Range("L136").Value = 214288290
$ws.Range("M136").Value = -9969690.600000001
$ws.Range("N136").Value = -214293390

$ws = $wb.Worksheets.Item("WVR")
# WVR!row25 - A Drag of a Doublet / Initiate's Doublet Vest
$ws.Range("H25").Value = 3100
$ws.Range("J25").Value = 3320
$ws.Range("L25").Value = 3320
$ws.Range("N25").Value = -3906

# WVR!row28 - Doublet Jeopardy / Cotton Doublet Vest of Gathering
$ws.Range("H28").Value = 3000
$ws.Range("J28").Value = 3000
$ws.Range("L28").Value = 3000
$ws.Range("N28").Value = -3696
